# Apply the WRESBAL.xlsx update:
#  - Data sheet: append new weekly observation row 95 (2023-07-19 -> 3230.457)
#  - SeriesInfo sheet: refresh the FRED metadata (realtime_start/end,
#    observation_end, last_updated, popularity)

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: add row 95 -------------------------------------------------
# Copy the formatting of the prior data row (date style with YYYY-MM-DD
# number format) onto the new row's date cell, then fill in the values.
$wsData.Range("A94").Copy()
$wsData.Range("A95").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsData.Cells.Item(95, 1).Value = 45126
$wsData.Cells.Item(95, 2).Value = 3230.457

# --- SeriesInfo sheet: update metadata values -------------------------------
# These cells hold plain text that merely looks like a date/timestamp (they
# are not real dates in the source file). Assigning a "YYYY-MM-DD" looking
# string straight to .Value would make Excel auto-convert it to a date
# serial number, which we don't want. Instead, write it as a text formula
# (="...") and then collapse the formula down to its static text result via
# copy / paste-values, which leaves the cell's formatting completely
# untouched (no number format / style changes).
function Set-TextValue($range, [string]$value) {
    $escaped = $value.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

Set-TextValue $wsInfo.Range("B3") "2023-07-24"
Set-TextValue $wsInfo.Range("B4") "2023-07-24"
Set-TextValue $wsInfo.Range("B7") "2023-07-19"
Set-TextValue $wsInfo.Range("B14") "2023-07-20 15:35:27-05"

$wsInfo.Range("B15").Value = 78
